$d = $word.ActiveDocument

# --- 1) "For longer probes..." note -> "Other supply voltages..." ---
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd() -eq "For longer probes please add `$ {{length_adder}} {{adder_per}}") {
        $para.Range.Text = "Other supply voltages available at no extra charge"
        $found = $true
        break
    }
}
if (-not $found) {
    throw "Could not find the 'For longer probes' paragraph"
}

# --- 2) Remove the (now duplicate) original "Other supply voltages..." note paragraph ---
# After step 1 there are two consecutive paragraphs with this text; delete the second one
# (the original bullet that used to carry this text), keeping list structure intact.
$targetIndex = -1
$seen = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd() -eq "Other supply voltages available at no extra charge") {
        $seen += 1
        if ($seen -eq 2) {
            $targetIndex = $i
            break
        }
    }
}
if ($targetIndex -eq -1) {
    throw "Could not find the duplicate 'Other supply voltages' paragraph"
}
$d.Paragraphs.Item($targetIndex).Range.Delete()

# --- 3) Merge "Delivery: " + "{{lead_time}}" + tabs/"FOB; Houston, TX" runs into one run ---
$deliveryPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("Delivery: ")) {
        $deliveryPara = $para
        break
    }
}
if ($null -eq $deliveryPara) {
    throw "Could not find the 'Delivery:' paragraph"
}

$xml = @"
<?xml version='1.0'?>
<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>
<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>
<pkg:xmlData>
<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
<w:body>
<w:p>
<w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:szCs w:val="24"/></w:rPr></w:pPr>
<w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>Delivery: {{lead_time}}</w:t><w:tab/><w:tab/><w:tab/><w:tab/><w:tab/><w:t>FOB; Houston, TX</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$deliveryPara.Range.InsertXML($xml)
